$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Day & Night shift dropdown"
$ws.Range("A5").Value = "Shared Services separation"

$ws.Range("C4").Value = "Phase 2"
$ws.Range("C5").Value = "Phase 2"
$ws.Range("C1").Value = "Phases"
$ws.Range("C2").Value = "Phase 1"
$ws.Range("C3").Value = "Phase 1"
